$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to Text format so numeric-looking strings
# (e.g. "1.94", "1.00") are preserved exactly as text, matching the source data.
$ws.Range("D5,D6,D7,D8,D9,D10,D12,D13,D14,D18,D20,D21,D22,D23,D24,D26,D27,D28,D29,D30,D31,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D45,D46,D47,D48,D50,D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range('D2').Value2 = '96.829.15'
$ws.Range('E2').Value2 = '  +0.77%  '
$ws.Range('D3').Value2 = '3.704.74'
$ws.Range('E3').Value2 = '  +4.49%  '
$ws.Range('E4').Value2 = '  +0.00%  '
$ws.Range('D5').Value2 = '244.39'
$ws.Range('E5').Value2 = '  +1.92%  '
$ws.Range('D6').Value2 = '1.94'
$ws.Range('E6').Value2 = '  +18.86%  '
$ws.Range('D7').Value2 = '673.81'
$ws.Range('E7').Value2 = '  +3.36%  '
$ws.Range('D8').Value2 = '0.428'
$ws.Range('E8').Value2 = '  +5.40%  '
$ws.Range('D9').Value2 = '1.13'
$ws.Range('E9').Value2 = '  +5.78%  '
$ws.Range('D10').Value2 = '1.00'
$ws.Range('E10').Value2 = '  -0.02%  '
$ws.Range('D11').Value2 = '3.700.60'
$ws.Range('E11').Value2 = '  +4.43%  '
$ws.Range('D12').Value2 = '45.41'
$ws.Range('E12').Value2 = '  +4.75%  '
$ws.Range('D13').Value2 = '0.206'
$ws.Range('E13').Value2 = '  +1.87%  '
$ws.Range('D14').Value2 = '6.60'
$ws.Range('E14').Value2 = '  +3.70%  '
$ws.Range('D15').Value2 = '4.393.56'
$ws.Range('E15').Value2 = '  +4.52%  '
$ws.Range('D16').Value2 = '96.605.86'
$ws.Range('E16').Value2 = '  +0.66%  '
$ws.Range('E17').Value2 = '  +2.03%  '
$ws.Range('D18').Value2 = '8.83'
$ws.Range('E18').Value2 = '  +12.19%  '
$ws.Range('D19').Value2 = '3.679.42'
$ws.Range('E19').Value2 = '  +4.76%  '
$ws.Range('D20').Value2 = '13.10'
$ws.Range('E20').Value2 = '  +5.90%  '
$ws.Range('D21').Value2 = '18.61'
$ws.Range('E21').Value2 = '  +5.67%  '
$ws.Range('D22').Value2 = '0.550'
$ws.Range('E22').Value2 = '  +3.93%  '
$ws.Range('D23').Value2 = '517.49'
$ws.Range('E23').Value2 = '  +2.52%  '
$ws.Range('D24').Value2 = '3.44'
$ws.Range('E24').Value2 = '  +1.74%  '
$ws.Range('E25').Value2 = '  +7.41%  '
$ws.Range('D26').Value2 = '6.92'
$ws.Range('E26').Value2 = '  +1.13%  '
$ws.Range('D27').Value2 = '101.67'
$ws.Range('E27').Value2 = '  +6.20%  '
$ws.Range('D28').Value2 = '13.05'
$ws.Range('E28').Value2 = '  +2.87%  '
$ws.Range('D29').Value2 = '0.170'
$ws.Range('E29').Value2 = '  +13.26%  '
$ws.Range('D30').Value2 = '3.09'
$ws.Range('E30').Value2 = '  +3.98%  '
$ws.Range('D31').Value2 = '12.13'
$ws.Range('E31').Value2 = '  +7.15%  '
$ws.Range('E32').Value2 = '  -0.15%  '
$ws.Range('D33').Value2 = '0.187'
$ws.Range('E33').Value2 = '  +2.98%  '
$ws.Range('D34').Value2 = '33.39'
$ws.Range('E34').Value2 = '  +7.19%  '
$ws.Range('D35').Value2 = '0.997'
$ws.Range('E35').Value2 = '  -0.21%  '
$ws.Range('D36').Value2 = '1.75'
$ws.Range('E36').Value2 = '  +9.58%  '
$ws.Range('D37').Value2 = '0.596'
$ws.Range('E37').Value2 = '  +6.18%  '
$ws.Range('B38').Value2 = 'RenderToken'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value2 = '8.83'
$ws.Range('E38').Value2 = '  +1.55%  '
$ws.Range('B39').Value2 = 'Bittensor'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value2 = '618.84'
$ws.Range('E39').Value2 = '  +1.34%  '
$ws.Range('D40').Value2 = '42.73'
$ws.Range('E40').Value2 = '  +28.81%  '
$ws.Range('D41').Value2 = '0.160'
$ws.Range('E41').Value2 = '  +7.09%  '
$ws.Range('D42').Value2 = '0.969'
$ws.Range('E42').Value2 = '  +8.26%  '
$ws.Range('E43').Value2 = '  +8.85%  '
$ws.Range('D45').Value2 = '6.13'
$ws.Range('E45').Value2 = '  +7.91%  '
$ws.Range('D46').Value2 = '0.0447'
$ws.Range('E46').Value2 = '  +6.35%  '
$ws.Range('D47').Value2 = '0.424'
$ws.Range('E47').Value2 = '  +25.02%  '
$ws.Range('D48').Value2 = '2.31'
$ws.Range('E48').Value2 = '  +2.12%  '
$ws.Range('E49').Value2 = '  +0.26%  '
$ws.Range('D50').Value2 = '8.61'
$ws.Range('E50').Value2 = '  +5.97%  '
$ws.Range('D51').Value2 = '54.62'
$ws.Range('E51').Value2 = '  +3.59%  '
